# Update "想去人数" (interested-attendee count) figures pulled from bilibili,
# refreshing the generated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 825
$wsExhibit.Range("F5").Value = 1016
$wsExhibit.Range("F6").Value = 2378

# Sheet "全部类型" (all types / combined view) mirrors the same rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 825
$wsAll.Range("F7").Value = 1016
$wsAll.Range("F8").Value = 2378
